$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.491.83'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.835.05'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'318.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = "'0.5319"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.50%  '
$ws.Range('D8').Value = "'0.4077"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.85%  '
$ws.Range('D9').Value = "'0.07568"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').Value = "'41.88"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').Value = "'1.111"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = "'6.333"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.96%  '
$ws.Range('D13').Value = "'7.633"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.46%  '
$ws.Range('D14').Value = "'1.001"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'20.85"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').Value = '1.838.34'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').Value = "'0.06596"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.75%  '
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').Value = "'1.000"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'6.074"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.57%  '
$ws.Range('D23').Value = '28.500.03'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = "'11.35"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').Value = "'2.113"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').Value = "'2.455"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.59%  '
$ws.Range('D27').Value = "'156.91"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.04%  '
$ws.Range('D28').Value = "'20.58"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('D29').Value = '2.048.17'
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('D30').Value = "'123.95"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('D31').Value = "'1.126"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.96%  '
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('D33').Value = "'5.695"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.58%  '
$ws.Range('D34').Value = "'3.658"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = "'0.07174"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.91%  '
$ws.Range('D36').Value = "'0.2273"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').Value = "'5.277"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.98%  '
$ws.Range('D38').Value = "'0.02350"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.11%  '
$ws.Range('D39').Value = "'8.855"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.59%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = "'11.36"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'0.6285"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.55%  '
$ws.Range('D42').Value = "'1.193"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.57%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = "'13.41"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').Value = "'3.718"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.96%  '
$ws.Range('D47').Value = "'0.5861"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('D48').Value = "'125.85"
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = "'1.991"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').Value = "'1.195"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').Value = "'0.06908"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.50%  '
